$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last (duplicate) row entirely, and shift data up by deleting the header row.
# Target layout:
# Row1: 192101031686 (number) | asd00asd! | Grund
# Row2: 198806011642 (number) | asd00asd! | Grund
# Row3: 198109231665 (number) | asd00asd! | Admin

# Delete row 5 (duplicate of row4) and row1 (header row)
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(1).Delete()

# Now rows are:
# Row1: 192101031686 | asd00asd! | Grund
# Row2: 198806011642 | asd00asd! | Grund
# Row3: 198109231665 | asd00asd! | Admin

# Ensure column A is numeric (not text) format, then set numeric values
$ws.Range("A1:A3").NumberFormat = "0"

$ws.Range("A1").Value = 192101031686
$ws.Range("A2").Value = 198806011642
$ws.Range("A3").Value = 198109231665

# Column A width now auto/best-fits the new numeric values
$ws.Columns.Item(1).ColumnWidth = 14.44140625

# Clear the selection set on A7 (select A1 instead, default state)
$ws.Range("A1").Select()
